# Update NATMI LR-pair TPM-derived metrics (Fn1-Tshr) with recomputed values
# from the new TPM expression matrix. Ligand/Receptor/Edge expression,
# specificity and weight columns (G,H,I,J,M,N,O,P,Q,R,S,T) are refreshed
# for rows 2-10 per the updated script output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 0.05057900000000001
$ws.Range("N2").Value = 0.151737
$ws.Range("O2").Value = 0.01400296657613869
$ws.Range("P2").Value = 0.01400296657613869
$ws.Range("Q2").Value = 0.315618422532
$ws.Range("R2").Value = 2.840565802788
$ws.Range("S2").Value = 0.0002425636612328911
$ws.Range("T2").Value = 0.0002425636612328911
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.146324388539341
$ws.Range("P3").Value = 0.146324388539341
$ws.Range("Q3").Value = 3.29806348088
$ws.Range("R3").Value = 29.68257132792
$ws.Range("S3").Value = 0.002534675721660817
$ws.Range("T3").Value = 0.002534675721660818
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("O4").Value = 0.8396726448845202
$ws.Range("P4").Value = 0.8396726448845202
$ws.Range("Q4").Value = 18.925715074784
$ws.Range("R4").Value = 170.331435673056
$ws.Range("S4").Value = 0.01454506585250005
$ws.Range("T4").Value = 0.01454506585250005
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("M5").Value = 0.05057900000000001
$ws.Range("N5").Value = 0.151737
$ws.Range("O5").Value = 0.01400296657613869
$ws.Range("P5").Value = 0.01400296657613869
$ws.Range("Q5").Value = 17.478411999241
$ws.Range("R5").Value = 157.305707993169
$ws.Range("S5").Value = 0.01343276343966564
$ws.Range("T5").Value = 0.01343276343966565
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.146324388539341
$ws.Range("P6").Value = 0.146324388539341
$ws.Range("S6").Value = 0.1403660350123243
$ws.Range("T6").Value = 0.1403660350123243
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("O7").Value = 0.8396726448845202
$ws.Range("P7").Value = 0.8396726448845202
$ws.Range("S7").Value = 0.8054810346196188
$ws.Range("T7").Value = 0.8054810346196191
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("M8").Value = 0.05057900000000001
$ws.Range("N8").Value = 0.151737
$ws.Range("O8").Value = 0.01400296657613869
$ws.Range("P8").Value = 0.01400296657613869
$ws.Range("Q8").Value = 0.4263171730213334
$ws.Range("R8").Value = 3.836854557192
$ws.Range("S8").Value = 0.0003276394752401566
$ws.Range("T8").Value = 0.0003276394752401567
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.146324388539341
$ws.Range("P9").Value = 0.146324388539341
$ws.Range("R9").Value = 40.09331836528
$ws.Range("S9").Value = 0.003423677805355895
$ws.Range("T9").Value = 0.003423677805355896
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("O10").Value = 0.8396726448845202
$ws.Range("P10").Value = 0.8396726448845202
$ws.Range("S10").Value = 0.01964654441240122
$ws.Range("T10").Value = 0.01964654441240122
